# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values are forced to remain text (matching the original inlineStr/text cell type)
# by prefixing with an apostrophe, then the style is reset to "Normal" so no stray
# number-format / quote-prefix styling is left behind on the cell.
$ws.Range("D2").Value = "'26.407.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.722.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'242.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4923"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.2615"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.06198"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'1.726.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07017"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'15.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'4.572"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.5994"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'77.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'26.404.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.9996"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007168"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'1.939.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'4.482"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'8.582"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'5.159"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'137.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'15.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'107.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'1.705"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Value = "'0.07963"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'3.665"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.04542"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'2.603"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.9941"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.6249"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'0.9237"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Value = "'1.946"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.9995"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.01483"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'99.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'5.325"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.3838"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'6.726"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.1164"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.05364"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'30.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'7.671"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'1.234"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'50.83"
$ws.Range("D51").Style = "Normal"

# Volume percentage values (already non-numeric text, e.g. "  -0.39%  ")
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("E39").Value = "  -5.92%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("E43").Value = "  -4.01%  "
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("E51").Value = "  -0.71%  "
